$wb = $excel.ActiveWorkbook

# --- "Fix" sheet: append a new row (row 4) ---------------------------------
$wsFix = $wb.Worksheets.Item("Fix")
$wsFix.Range("A4").Value = 202307
$wsFix.Range("B4").Value = 133.75
$wsFix.Range("C4").Value = 10.85
$wsFix.Range("D4").Value = 145.22

# --- "Flow" sheet: new "Remark" column (G) + a new data row (10) -----------
$wsFlow = $wb.Worksheets.Item("Flow")

# Write the new shared strings in the same first-seen order as the source
# workbook so the shared-string table comes out in the same order.
$wsFlow.Range("G1").Value = "Remark"

# Match the highlighted-header formatting used by the rest of row 1.
$wsFlow.Range("F1").Copy()
$wsFlow.Range("G1").PasteSpecial(-4122)

$wsFlow.Range("G9").Value = "Indah Water"
$wsFlow.Range("G6").Value = "Babadona"
$wsFlow.Range("G5").Value = "Kitchen sink"
$wsFlow.Range("G2").Value = "Medeena"

$wsFlow.Range("G3").Value = "Medeena"
$wsFlow.Range("G4").Value = "Medeena"
$wsFlow.Range("G7").Value = "Babadona"
$wsFlow.Range("G8").Value = "Babadona"

# New row 10
$wsFlow.Range("A10").Value = 202307
$wsFlow.Range("B10").Value = "Amir"
$wsFlow.Range("C10").Value = "Others"
$wsFlow.Range("D10").Value = 60
$wsFlow.Range("E10").Value = 5
$wsFlow.Range("F10").Value = "Adnan, Amir, Kimi, Lutfi, Ziad"
$wsFlow.Range("G10").Value = "Indah Water"

# --- sheet view / active-tab bookkeeping ------------------------------------
# Active tab moves from "Rent" to "Fix"; Fix's selection lands on N5.
$null = $wsFix.Range("N5").Select()

# Flow's selection moves too (no longer the active tab, but Excel still
# tracks a remembered selection per sheet).
$null = $wsFlow.Range("H19").Select()

# Finally, re-select back on Fix so it remains the active/visible tab
# (matches the saved file: Fix is tabSelected, Rent is not).
$wsFix.Select()
